$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..11, columns B:G (ME, MAE, MSE, RMSE, SE, N)
# Each row shifts up by one (row N gets what used to be row N+1),
# and the last row (11) gets a brand-new data point appended at the
# end of the ifoCAST evaluation series.
$data = @(
    @(-0.07072591166527165, 0.5579505877006946, 0.6529064599314954, 0.8080262742828945, 0.8282610236588537, 18),
    @(0.03916263344646955, 0.552040872335158, 0.6357226611696118, 0.7973221815361792, 0.8208689046632819, 17),
    @(0.05955592157360488, 0.532271854806013, 0.7198117377579547, 0.8484171955812509, 0.8740799841134135, 16),
    @(0.09235442246305205, 0.5824508295849771, 0.7506227629998017, 0.8663848815623468, 0.8916838547205332, 15),
    @(0.1521483464337763, 0.6635260061488867, 1.004366658952618, 1.002180951202236, 1.0279571060649, 14),
    @(0.159095420907867, 0.6746913222580827, 0.8632158529041898, 0.9290941033631576, 0.9527485895797485, 13),
    @(0.2207410413997515, 0.686169913841685, 0.7776458313225046, 0.8818422939066285, 0.8917312420012032, 12),
    @(0.220829446489595, 0.7371449325608843, 1.012670326230143, 1.006315222099985, 1.029706341620059, 11),
    @(0.1936765229212045, 0.6934653330366425, 1.055080256735366, 1.027170996833227, 1.063312257473744, 10),
    @(0.2439803621284286, 0.799271525189887, 1.113210894331233, 1.05508809790047, 1.088758484182915, 9)
)

$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $vals[$j]
    }
}
